$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 137.5335740648586
    3  = 12.74799120842664
    4  = 13.04650845981588
    5  = 17.66239860661856
    6  = 22.65584784772332
    7  = 7.489790163390561
    8  = 8.487405656465125
    9  = 22.62975636559675
    10 = 38.99295264161898
    11 = 10.86249792651572
    12 = 2.932222153106079
    13 = 6.631073442813303
    14 = 1.941513228829725
    15 = 3.026611926681596
    16 = 19.38904080617069
    17 = 19.97763394708468
    18 = 19.62616633490918
    19 = 6.267327486107655
    20 = 24.84906713942108
    21 = 70.10167330181636
    22 = 10.93463320063034
    23 = 2.331350667661939
    24 = 22.56222547068093
    25 = 6.839037903292694
    26 = 13.02271916728871
    27 = 23.91054117939786
    28 = 5.345684249813386
    29 = 11.27535726134196
    30 = 2.682818279837439
    31 = 2.672074728373559
    32 = 4.906733432860579
    33 = 5.123906651737579
    34 = 93.66995562743634
    35 = 7.914927842746952
    36 = 22.84385999834121
    37 = 3.942883387243925
    38 = 10.04215103259517
    39 = 9.750540350004147
    40 = 7.623317160155925
    41 = 5.95422970058887
    42 = 259.8
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 3).Value = $newValues[$row]
}
